# Hide "note" type fields from the contents screen:
# add a new `hideInContents` column to the survey sheet and flag the
# existing note row (row 6 - the "note" field) as hidden.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# New header cell for the hideInContents column (column H).
$ws.Range("H1").Value = "hideInContents"

# Flag the note row (A6 = "note") to be hidden from the contents screen.
$ws.Range("H6").Value = $true

# The survey sheet becomes the active sheet/selection (it was "settings").
$ws.Activate() | Out-Null
$ws.Range("H7").Select() | Out-Null
